# Ele.xlsx edit: bump the "用电 KWh" (electricity usage) column C values
# from 100 to their new readings (row 3 -> 3000, all the rest -> 1000),
# then leave the selection on C4:C11 (active cell C4) as the user left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1000
    3  = 3000
    4  = 1000
    5  = 1000
    6  = 1000
    7  = 1000
    8  = 1000
    9  = 1000
    10 = 1000
    11 = 1000
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

$ws.Range("C4:C11").Select()
